# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 299 of Sheet1, shifting all
# subsequent rows (old 299-380) down by one (new 300-381).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(299).Insert()

$ws.Range("A299").Value = 11
$ws.Range("B299").Value = "Vega Monumental Concepción"
$ws.Range("C299").Value = "Bíobío"
$ws.Range("D299").Value = 44524
$ws.Range("E299").Value = 8
$ws.Range("F299").Value = "Fruta"
$ws.Range("G299").Value = 100106
$ws.Range("H299").Value = "Oleaginosos"
$ws.Range("I299").Value = 100106002
$ws.Range("J299").Value = "Palta"
$ws.Range("K299").Value = "Hass"
$ws.Range("L299").Value = "Segunda"
$ws.Range("M299").Value = 180
$ws.Range("N299").Value = 2000
$ws.Range("O299").Value = 2200
$ws.Range("P299").Value = 2111
$ws.Range("Q299").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R299").Value = "Región de O'Higgins"
$ws.Range("S299").Value = 2111
$ws.Range("T299").Value = 1
